$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value() = "'59.367.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value() = "  +0.72%  "
$ws.Range("D3").Value() = "'2.601.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value() = "  +0.68%  "
$ws.Range("D5").Value() = "'535.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value() = "  +2.81%  "
$ws.Range("D6").Value() = "'141.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value() = "  +1.89%  "
$ws.Range("E7").Value() = "  +0.11%  "
$ws.Range("D8").Value() = "'0.566"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value() = "  +0.59%  "
$ws.Range("D9").Value() = "'6.49"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value() = "  -1.14%  "
$ws.Range("E10").Value() = "  +1.51%  "
$ws.Range("E11").Value() = "  +1.67%  "
$ws.Range("E12").Value() = "  -0.67%  "
$ws.Range("D13").Value() = "'3.064.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value() = "  +0.72%  "
$ws.Range("D14").Value() = "'59.289.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value() = "  +0.66%  "
$ws.Range("D15").Value() = "'20.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value() = "  +1.16%  "
$ws.Range("D16").Value() = "'2.658.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value() = "  +3.34%  "
$ws.Range("E17").Value() = "  +0.45%  "
$ws.Range("D18").Value() = "'341.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value() = "  +1.19%  "
$ws.Range("E19").Value() = "  +1.67%  "
$ws.Range("D20").Value() = "'10.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value() = "  +0.09%  "
$ws.Range("E21").Value() = "  -2.26%  "
$ws.Range("D22").Value() = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value() = "  +0.05%  "
$ws.Range("D23").Value() = "'67.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value() = "  +2.40%  "
$ws.Range("E24").Value() = "  +1.62%  "
$ws.Range("E25").Value() = "  -1.43%  "
$ws.Range("D26").Value() = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value() = "  +0.10%  "
$ws.Range("D27").Value() = "'7.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value() = "  +3.19%  "
$ws.Range("D28").Value() = "'0.0₃0743"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value() = "  +2.72%  "
$ws.Range("E29").Value() = "  +0.04%  "
$ws.Range("E30").Value() = "  +5.94%  "
$ws.Range("E31").Value() = "  -1.74%  "
$ws.Range("E32").Value() = "  +0.86%  "
$ws.Range("D33").Value() = "'150.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value() = "  +0.67%  "
$ws.Range("E34").Value() = "  +0.31%  "
$ws.Range("E35").Value() = "  -0.62%  "
$ws.Range("E36").Value() = "  -0.29%  "
$ws.Range("E37").Value() = "  +3.58%  "
$ws.Range("D38").Value() = "'0.825"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value() = "  +0.44%  "
$ws.Range("D39").Value() = "'3.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value() = "  +0.93%  "
$ws.Range("E40").Value() = "  +0.14%  "
$ws.Range("D41").Value() = "'273.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value() = "  +1.32%  "
$ws.Range("D42").Value() = "'0.598"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value() = "  +1.83%  "
$ws.Range("D43").Value() = "'10.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value() = "  -0.15%  "
$ws.Range("D44").Value() = "'0.0953"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value() = "  +0.12%  "
$ws.Range("E45").Value() = "  +1.27%  "
$ws.Range("B46").Value() = "'InjectiveProtocol"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value() = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value() = "'18.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value() = "  +3.56%  "
$ws.Range("B47").Value() = "'Maker"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value() = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value() = "'1.948.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value() = "  -0.68%  "
$ws.Range("E48").Value() = "  +1.48%  "
$ws.Range("D49").Value() = "'4.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value() = "  +0.10%  "
$ws.Range("D50").Value() = "'111.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value() = "  -2.25%  "
$ws.Range("D51").Value() = "'4.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value() = "  +0.39%  "
